# FEAT: Optimizacion para las columnas de las interfaces
#
# The sheet used to keep one "interfaces" column holding a stringified
# Python list-of-dicts per device. This splits that blob into three
# columns (IP / STATUS / PROTO) for each of the four interfaces that
# actually carry data: GigabitEthernet0/0, GigabitEthernet0/1,
# FastEthernet0/2/0 and FastEthernet0/2/1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header row ----------------------------------------------------
# I1 used to read "interfaces"; it becomes the first of the new per-
# interface headers, and J1:T1 are brand new columns.
$ws.Range("I1").Value = "GigabitEthernet0/0_IP"
$ws.Range("J1").Value = "GigabitEthernet0/0_STATUS"
$ws.Range("K1").Value = "GigabitEthernet0/0_PROTO"
$ws.Range("L1").Value = "GigabitEthernet0/1_IP"
$ws.Range("M1").Value = "GigabitEthernet0/1_STATUS"
$ws.Range("N1").Value = "GigabitEthernet0/1_PROTO"
$ws.Range("O1").Value = "FastEthernet0/2/0_IP"
$ws.Range("P1").Value = "FastEthernet0/2/0_STATUS"
$ws.Range("Q1").Value = "FastEthernet0/2/0_PROTO"
$ws.Range("R1").Value = "FastEthernet0/2/1_IP"
$ws.Range("S1").Value = "FastEthernet0/2/1_STATUS"
$ws.Range("T1").Value = "FastEthernet0/2/1_PROTO"

# Carry the existing header formatting (bold font, thin border, centered/
# top-aligned) from I1 onto the newly-added header cells.
$ws.Range("I1").Copy()
$ws.Range("J1:T1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: first device ---------------------------------------------------
# Domain updated and the old JSON-ish "interfaces" blob is replaced by its
# parsed-out per-interface values.
$ws.Range("H2").Value = "jane.com"

$ws.Range("I2").Value = "192.168.2.1"
$ws.Range("J2").Value = "administratively down"
$ws.Range("K2").Value = "down"
$ws.Range("L2").Value = "unassigned"
$ws.Range("M2").Value = "administratively down"
$ws.Range("N2").Value = "down"
$ws.Range("O2").Value = "unassigned"
$ws.Range("P2").Value = "administratively down"
$ws.Range("Q2").Value = "down"
$ws.Range("R2").Value = "unassigned"
$ws.Range("S2").Value = "administratively down"
$ws.Range("T2").Value = "down"

# --- Row 3: second device ---------------------------------------------------
# I3 was already present as an empty placeholder cell; replicate that same
# blank placeholder across the new per-interface columns J3:T3.
$ws.Range("I3").Copy($ws.Range("J3:T3"))
$excel.CutCopyMode = $false
